# Update the workbook:
# 1) Column C ("Förändrad") date serial value changes from 45184 to 45186 for every
#    data row.
# 2) The HYPERLINK() formulas in columns S-Y get a second argument added: the
#    "friendly name" shown for the link, which is equal to the row's "Beteckning"
#    value in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = 2; $row -le $lastRow; $row++) {

    $cCell = $ws.Range("C" + $row)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    $name = $ws.Range("A" + $row).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $row)
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            if ($formula.EndsWith(")") -and -not $formula.Contains(",")) {
                $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $name + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}

Write-Output "done"
